$wb = $excel.ActiveWorkbook

# --- LinkedList sheet: add the 4 new "type of linked list" rows ---
$ws = $wb.Worksheets.Item("LinkedList")
$ws.Cells.Item(10, 1).Value = "singly linked list"
$ws.Cells.Item(11, 1).Value = "linked list"
$ws.Cells.Item(12, 1).Value = "doubly linked list"
$ws.Cells.Item(13, 1).Value = "type of the linked list"

# Make LinkedList the active sheet/tab, and select A14 (next empty row)
$ws.Activate()
$ws.Range("A14").Select()
